$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scores for "List User Ads" (row 21) and "Logout" (row 31)
$ws.Range("C21").Value = 2
$ws.Range("C31").Value = 5

# Recalculate so the cached SUM formula value (C51) reflects the change
$excel.Calculate()

# Move the active selection to C21, matching the saved view state
$ws.Range("C21").Select()
